$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.972.08"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.218.78"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.01"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.02"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.45"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.35"
$ws.Range("E11").Value = "  +5.96%  "
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("E13").Value = "  +3.56%  "
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "2.559.19"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.76"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "2.220.96"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.732"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "39.885.28"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.15"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.75"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.56"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.43"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.12"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.23"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("E30").Value = "  -7.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.02"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.84"
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("E35").Value = "  +6.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0715"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0992"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.73"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.21"
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("D42").Value = "2.092.14"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.96"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.80"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.00"
$ws.Range("E47").Value = "  -7.94%  "
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D49").Value = "2.433.40"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  +2.66%  "
